# Update the cached text of the "datetimeFigureOut" date/time fields that
# live on the slide master, every slide layout, and the notes master.
#
# The source deck was re-saved on a different day (2020-12-26 -> 2021-07-25),
# and PowerPoint recached the auto date field's displayed text accordingly:
#   - English (en-US) placeholders:  12/26/2020  -> 7/25/2021
#   - Korean  (ko-KR) placeholder:   2020-12-26  -> 2021-07-25
#
# ppPlaceholderDate = 16

$ppPlaceholderDate = 16
$usDate = "7/25/2021"
$krDate = "2021-07-25"

function Update-DatePlaceholder {
    param($shapes)

    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame -eq -1 -and $shp.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
            $current = $shp.TextFrame.TextRange.Text
            if ($current -eq "2020-12-26") {
                $shp.TextFrame.TextRange.Text = $krDate
            } else {
                $shp.TextFrame.TextRange.Text = $usDate
            }
        }
    }
}

$p = $ppt.ActivePresentation

# Slide master
Update-DatePlaceholder $p.SlideMaster.Shapes

# Every slide layout under the master
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DatePlaceholder $layouts.Item($li).Shapes
}

# Notes master (Korean locale field)
Update-DatePlaceholder $p.NotesMaster.Shapes
